$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.205.30"
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("D3").Value = "2.244.22"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.53"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.625"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.39%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "74.31"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.67%  "
$ws.Range("E8").Value = "  +0.04%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.617"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -3.43%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "42.08"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.31%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0944"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.11"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.83%  "
$ws.Range("E13").Value = "  -2.13%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "14.48"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.51%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.854"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.24%  "
$ws.Range("D16").Value = "2.266.07"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "42.063.48"
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("D18").Value = "0.0₃0987"
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("E19").Value = "  -0.48%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.24"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "230.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.62"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +33.67%  "
$ws.Range("E24").Value = "  +0.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.41"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.61"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -4.20%  "
$ws.Range("E27").Value = "  -3.11%  "
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "169.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.08"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -5.02%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "20.64"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0820"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.27%  "
$ws.Range("E32").Value = "  -3.52%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.48"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.25%  "
$ws.Range("E34").Value = "  -1.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +9.71%  "
$ws.Range("E36").Value = "  -0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0307"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.47%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "13.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.73%  "
$ws.Range("E39").Value = "  -3.17%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "62.24"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.53%  "
$ws.Range("E42").Value = "  -2.06%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "107.46"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.51%  "
$ws.Range("E44").Value = "  +2.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.65"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.997"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.14%  "
$ws.Range("B47").Value = "FTXToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.40%  "
$ws.Range("B48").Value = "ARBITRUM"
$ws.Range("C48").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.12"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -2.34%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("E50").Value = "  +1.49%  "
$ws.Range("E51").Value = "  +0.10%  "
